$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), matching the header style used by
# the rest of row 1 (bold/bordered/centered -> same style as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-17 for the two new columns, I (I0) and J (IF).
$iValues = @(6, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 5)
$jValues = @(8, 4, 3, 5, 6, 5, 5, 3, 6, 5, 6, 6, 7, 5, 3, 5)

for ($r = 2; $r -le 17; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
